$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 60,6

$data[0,0] = 'Response group type'
$data[0,1] = 'ratio_paras'
$data[0,2] = 'One-way F test'
$data[0,3] = [double]'103.9250343758279'
$data[0,4] = [double]'5.450084891987035e-43'
$data[0,5] = $true

$data[1,0] = 'Response group type'
$data[1,1] = 'ratio_list_items'
$data[1,2] = 'One-way F test'
$data[1,3] = [double]'492.2975528874661'
$data[1,4] = [double]'4.982197960906943e-165'
$data[1,5] = $true

$data[2,0] = 'Response group type'
$data[2,1] = 'ratio_headings'
$data[2,2] = 'One-way F test'
$data[2,3] = [double]'2488.279675972683'
$data[2,4] = 0
$data[2,5] = $true

$data[3,0] = 'Response group type'
$data[3,1] = 'avg_para_len'
$data[3,2] = 'One-way F test'
$data[3,3] = [double]'36.80602324941301'
$data[3,4] = [double]'2.488443800449051e-16'
$data[3,5] = $true

$data[4,0] = 'Response group type'
$data[4,1] = 'num_sentences'
$data[4,2] = 'One-way F test'
$data[4,3] = [double]'1865.500228308488'
$data[4,4] = 0
$data[4,5] = $true

$data[5,0] = 'Response group type'
$data[5,1] = 'avg_len'
$data[5,2] = 'One-way F test'
$data[5,3] = [double]'340.6284512531364'
$data[5,4] = [double]'1.199104472394488e-122'
$data[5,5] = $true

$data[6,0] = 'Response group type'
$data[6,1] = 'flesch'
$data[6,2] = 'One-way F test'
$data[6,3] = [double]'93.1453954860687'
$data[6,4] = [double]'7.471291394604005e-39'
$data[6,5] = $true

$data[7,0] = 'Response group type'
$data[7,1] = 'cli'
$data[7,2] = 'One-way F test'
$data[7,3] = [double]'17.23422415869616'
$data[7,4] = [double]'3.982348758822768e-08'
$data[7,5] = $true

$data[8,0] = 'Response group type'
$data[8,1] = 'avg_concrete'
$data[8,2] = 'One-way F test'
$data[8,3] = [double]'4.10521155535292'
$data[8,4] = [double]'0.01667249606816113'
$data[8,5] = $true

$data[9,0] = 'Response group type'
$data[9,1] = 'concrete_ratio'
$data[9,2] = 'One-way F test'
$data[9,3] = [double]'24.28869629620745'
$data[9,4] = [double]'4.159931931841053e-11'
$data[9,5] = $true

$data[10,0] = 'Response group type'
$data[10,1] = 'abstract_ratio'
$data[10,2] = 'One-way F test'
$data[10,3] = [double]'259.0150279682046'
$data[10,4] = [double]'2.504918173677693e-97'
$data[10,5] = $true

$data[11,0] = 'Response group type'
$data[11,1] = 'undefined_ratio'
$data[11,2] = 'One-way F test'
$data[11,3] = [double]'292.7912155427557'
$data[11,4] = [double]'4.783773674273814e-108'
$data[11,5] = $true

$data[12,0] = 'SERP'
$data[12,1] = 'ratio_paras'
$data[12,2] = 'T-Test'
$data[12,3] = [double]'-3.20007098844747'
$data[12,4] = [double]'0.001461579957827798'
$data[12,5] = $true

$data[13,0] = 'SERP'
$data[13,1] = 'ratio_list_items'
$data[13,2] = 'T-Test'
$data[13,3] = [double]'-4.62777566626418'
$data[13,4] = [double]'4.72011467708349e-06'
$data[13,5] = $true

$data[14,0] = 'SERP'
$data[14,1] = 'ratio_headings'
$data[14,2] = 'T-Test'
$data[14,3] = [double]'5.292175244158653'
$data[14,4] = [double]'1.815062139272068e-07'
$data[14,5] = $true

$data[15,0] = 'SERP'
$data[15,1] = 'avg_para_len'
$data[15,2] = 'T-Test'
$data[15,3] = [double]'-22.93405871035168'
$data[15,4] = [double]'5.564501848320325e-80'
$data[15,5] = $true

$data[16,0] = 'SERP'
$data[16,1] = 'num_sentences'
$data[16,2] = 'T-Test'
$data[16,3] = [double]'-4.366461921739277'
$data[16,4] = [double]'1.536775195372723e-05'
$data[16,5] = $true

$data[17,0] = 'SERP'
$data[17,1] = 'avg_len'
$data[17,2] = 'T-Test'
$data[17,3] = [double]'-20.99987912806521'
$data[17,4] = [double]'1.358955706615782e-70'
$data[17,5] = $true

$data[18,0] = 'SERP'
$data[18,1] = 'flesch'
$data[18,2] = 'T-Test'
$data[18,3] = [double]'2.897497209476683'
$data[18,4] = [double]'0.003927216464984266'
$data[18,5] = $true

$data[19,0] = 'SERP'
$data[19,1] = 'cli'
$data[19,2] = 'T-Test'
$data[19,3] = [double]'-3.703154098703998'
$data[19,4] = [double]'0.0002367410099330563'
$data[19,5] = $true

$data[20,0] = 'SERP'
$data[20,1] = 'avg_concrete'
$data[20,2] = 'T-Test'
$data[20,3] = [double]'3.557182660562407'
$data[20,4] = [double]'0.000410682646011469'
$data[20,5] = $true

$data[21,0] = 'SERP'
$data[21,1] = 'concrete_ratio'
$data[21,2] = 'T-Test'
$data[21,3] = [double]'-0.867142270016976'
$data[21,4] = [double]'0.3862816838368289'
$data[21,5] = $false

$data[22,0] = 'SERP'
$data[22,1] = 'abstract_ratio'
$data[22,2] = 'T-Test'
$data[22,3] = [double]'-15.17561716272611'
$data[22,4] = [double]'4.968122309413295e-43'
$data[22,5] = $true

$data[23,0] = 'SERP'
$data[23,1] = 'undefined_ratio'
$data[23,2] = 'T-Test'
$data[23,3] = [double]'15.41329739875774'
$data[23,4] = [double]'4.144656743666704e-44'
$data[23,5] = $true

$data[24,0] = 'RR'
$data[24,1] = 'ratio_paras'
$data[24,2] = 'T-Test'
$data[24,3] = [double]'3.63490470761788'
$data[24,4] = [double]'0.0003069922751638511'
$data[24,5] = $true

$data[25,0] = 'RR'
$data[25,1] = 'ratio_list_items'
$data[25,2] = 'T-Test'
$data[25,3] = [double]'-3.879525504622102'
$data[25,4] = [double]'0.0001187192107984613'
$data[25,5] = $true

$data[26,0] = 'RR'
$data[26,1] = 'ratio_headings'
$data[26,2] = 'T-Test'
$data[26,3] = [double]'3.703629135984921'
$data[26,4] = [double]'0.000236309863528453'
$data[26,5] = $true

$data[27,0] = 'RR'
$data[27,1] = 'avg_para_len'
$data[27,2] = 'T-Test'
$data[27,3] = [double]'0.5409203513035115'
$data[27,4] = [double]'0.5888043460659246'
$data[27,5] = $false

$data[28,0] = 'RR'
$data[28,1] = 'num_sentences'
$data[28,2] = 'T-Test'
$data[28,3] = [double]'-3.226509327587789'
$data[28,4] = [double]'0.001335558424313903'
$data[28,5] = $true

$data[29,0] = 'RR'
$data[29,1] = 'avg_len'
$data[29,2] = 'T-Test'
$data[29,3] = [double]'3.789572000834149'
$data[29,4] = [double]'0.0001693798411145114'
$data[29,5] = $true

$data[30,0] = 'RR'
$data[30,1] = 'flesch'
$data[30,2] = 'T-Test'
$data[30,3] = [double]'-0.2054382678130561'
$data[30,4] = [double]'0.8373136719515815'
$data[30,5] = $false

$data[31,0] = 'RR'
$data[31,1] = 'cli'
$data[31,2] = 'T-Test'
$data[31,3] = [double]'0.4464872046650555'
$data[31,4] = [double]'0.6554394366694096'
$data[31,5] = $false

$data[32,0] = 'RR'
$data[32,1] = 'avg_concrete'
$data[32,2] = 'T-Test'
$data[32,3] = [double]'0.6760445824019018'
$data[32,4] = [double]'0.4993261673274386'
$data[32,5] = $false

$data[33,0] = 'RR'
$data[33,1] = 'concrete_ratio'
$data[33,2] = 'T-Test'
$data[33,3] = [double]'1.2416706355686'
$data[33,4] = [double]'0.2149425882217344'
$data[33,5] = $false

$data[34,0] = 'RR'
$data[34,1] = 'abstract_ratio'
$data[34,2] = 'T-Test'
$data[34,3] = [double]'1.967883838387863'
$data[34,4] = [double]'0.04963582312876014'
$data[34,5] = $true

$data[35,0] = 'RR'
$data[35,1] = 'undefined_ratio'
$data[35,2] = 'T-Test'
$data[35,3] = [double]'-1.210112196764809'
$data[35,4] = [double]'0.22680991195355'
$data[35,5] = $false

$data[36,0] = 'Chatbot'
$data[36,1] = 'ratio_paras'
$data[36,2] = 'T-Test'
$data[36,3] = [double]'-26.70903432385005'
$data[36,4] = [double]'3.478876148824275e-98'
$data[36,5] = $true

$data[37,0] = 'Chatbot'
$data[37,1] = 'ratio_list_items'
$data[37,2] = 'T-Test'
$data[37,3] = [double]'22.40540323826895'
$data[37,4] = [double]'2.050111511544678e-77'
$data[37,5] = $true

$data[38,0] = 'Chatbot'
$data[38,1] = 'ratio_headings'
$data[38,2] = 'T-Test'
$data[38,3] = [double]'18.07577900150219'
$data[38,4] = [double]'1.587913225899474e-56'
$data[38,5] = $true

$data[39,0] = 'Chatbot'
$data[39,1] = 'avg_para_len'
$data[39,2] = 'T-Test'
$data[39,3] = [double]'-20.3533957131679'
$data[39,4] = [double]'1.836427519341772e-67'
$data[39,5] = $true

$data[40,0] = 'Chatbot'
$data[40,1] = 'num_sentences'
$data[40,2] = 'T-Test'
$data[40,3] = [double]'23.35243727784407'
$data[40,4] = [double]'5.19572498798019e-82'
$data[40,5] = $true

$data[41,0] = 'Chatbot'
$data[41,1] = 'avg_len'
$data[41,2] = 'T-Test'
$data[41,3] = [double]'-24.28682441786492'
$data[41,4] = [double]'1.547847087079104e-86'
$data[41,5] = $true

$data[42,0] = 'Chatbot'
$data[42,1] = 'flesch'
$data[42,2] = 'T-Test'
$data[42,3] = [double]'-2.016040750719476'
$data[42,4] = [double]'0.04433218659223444'
$data[42,5] = $true

$data[43,0] = 'Chatbot'
$data[43,1] = 'cli'
$data[43,2] = 'T-Test'
$data[43,3] = [double]'2.659495576454042'
$data[43,4] = [double]'0.008077905587840311'
$data[43,5] = $true

$data[44,0] = 'Chatbot'
$data[44,1] = 'avg_concrete'
$data[44,2] = 'T-Test'
$data[44,3] = [double]'-0.8638296036201208'
$data[44,4] = [double]'0.3880973580777607'
$data[44,5] = $false

$data[45,0] = 'Chatbot'
$data[45,1] = 'concrete_ratio'
$data[45,2] = 'T-Test'
$data[45,3] = [double]'-1.921300344464408'
$data[45,4] = [double]'0.05526432420963388'
$data[45,5] = $false

$data[46,0] = 'Chatbot'
$data[46,1] = 'abstract_ratio'
$data[46,2] = 'T-Test'
$data[46,3] = [double]'-13.99724859110737'
$data[46,4] = [double]'8.894914808069599e-38'
$data[46,5] = $true

$data[47,0] = 'Chatbot'
$data[47,1] = 'undefined_ratio'
$data[47,2] = 'T-Test'
$data[47,3] = [double]'13.79218912698009'
$data[47,4] = [double]'7.004891698637221e-37'
$data[47,5] = $true

$data[48,0] = 'IAS'
$data[48,1] = 'ratio_paras'
$data[48,2] = 'One-way F test'
$data[48,3] = [double]'337.1931345863907'
$data[48,4] = [double]'5.072061204772169e-242'
$data[48,5] = $true

$data[49,0] = 'IAS'
$data[49,1] = 'ratio_list_items'
$data[49,2] = 'One-way F test'
$data[49,3] = [double]'503.0133287261479'
$data[49,4] = [double]'4.474863e-317'
$data[49,5] = $true

$data[50,0] = 'IAS'
$data[50,1] = 'ratio_headings'
$data[50,2] = 'One-way F test'
$data[50,3] = [double]'1434.350177387523'
$data[50,4] = 0
$data[50,5] = $true

$data[51,0] = 'IAS'
$data[51,1] = 'avg_para_len'
$data[51,2] = 'One-way F test'
$data[51,3] = [double]'121.9227978430191'
$data[51,4] = [double]'2.327858259410091e-108'
$data[51,5] = $true

$data[52,0] = 'IAS'
$data[52,1] = 'num_sentences'
$data[52,2] = 'One-way F test'
$data[52,3] = [double]'854.5093066193562'
$data[52,4] = 0
$data[52,5] = $true

$data[53,0] = 'IAS'
$data[53,1] = 'avg_len'
$data[53,2] = 'One-way F test'
$data[53,3] = [double]'600.8964257179297'
$data[53,4] = 0
$data[53,5] = $true

$data[54,0] = 'IAS'
$data[54,1] = 'flesch'
$data[54,2] = 'One-way F test'
$data[54,3] = [double]'40.10816492550137'
$data[54,4] = [double]'8.742930203466194e-39'
$data[54,5] = $true

$data[55,0] = 'IAS'
$data[55,1] = 'cli'
$data[55,2] = 'One-way F test'
$data[55,3] = [double]'11.30245280737467'
$data[55,4] = [double]'1.017809412477035e-10'
$data[55,5] = $true

$data[56,0] = 'IAS'
$data[56,1] = 'avg_concrete'
$data[56,2] = 'One-way F test'
$data[56,3] = [double]'2.65397987484021'
$data[56,4] = [double]'0.02136283399874483'
$data[56,5] = $true

$data[57,0] = 'IAS'
$data[57,1] = 'concrete_ratio'
$data[57,2] = 'One-way F test'
$data[57,3] = [double]'11.670469474077'
$data[57,4] = [double]'4.38810462932117e-11'
$data[57,5] = $true

$data[58,0] = 'IAS'
$data[58,1] = 'abstract_ratio'
$data[58,2] = 'One-way F test'
$data[58,3] = [double]'221.1192272393743'
$data[58,4] = [double]'8.646651848982935e-177'
$data[58,5] = $true

$data[59,0] = 'IAS'
$data[59,1] = 'undefined_ratio'
$data[59,2] = 'One-way F test'
$data[59,3] = [double]'228.4717117168348'
$data[59,4] = [double]'2.473556392891394e-181'
$data[59,5] = $true

$ws.Range("A2:F61").Value = $data

Write-Output "Updated rows 2:61"
